$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "584.17") must be
# forced to stay as text, matching the original inlineStr/text cell content,
# otherwise Excel auto-converts them to floating point numbers.
$textCells = @(
    'D5',
    'D6',
    'D10',
    'D14',
    'D18',
    'D19',
    'D20',
    'D21',
    'D23',
    'D24',
    'D26',
    'D27',
    'D29',
    'D30',
    'D31',
    'D32',
    'D36',
    'D38',
    'D39',
    'D40',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D51',
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '67.911.90'
$ws.Range('E2').Value = '  +1.35%  '
$ws.Range('D3').Value = '3.257.14'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '584.17'
$ws.Range('E5').Value = '  +1.01%  '
$ws.Range('D6').Value = '184.53'
$ws.Range('E6').Value = '  +4.32%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.92%  '
$ws.Range('D10').Value = '6.68'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('E11').Value = '  +1.87%  '
$ws.Range('D12').Value = '3.820.78'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('D14').Value = '28.66'
$ws.Range('E14').Value = '  +2.82%  '
$ws.Range('D15').Value = '67.918.01'
$ws.Range('E15').Value = '  +1.35%  '
$ws.Range('E16').Value = '  +2.82%  '
$ws.Range('D17').Value = '3.254.48'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '5.85'
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('D19').Value = '13.60'
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('D20').Value = '382.46'
$ws.Range('E20').Value = '  +3.40%  '
$ws.Range('D21').Value = '7.67'
$ws.Range('E21').Value = '  +1.13%  '
$ws.Range('E22').Value = '  -0.59%  '
$ws.Range('D23').Value = '71.36'
$ws.Range('E23').Value = '  +0.76%  '
$ws.Range('D24').Value = '0.514'
$ws.Range('E24').Value = '  +1.27%  '
$ws.Range('E25').Value = '  +0.96%  '
$ws.Range('D26').Value = '9.83'
$ws.Range('D27').Value = '0.182'
$ws.Range('E27').Value = '  +1.62%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').Value = '1.99'
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('D30').Value = '5.69'
$ws.Range('E30').Value = '  +1.16%  '
$ws.Range('D31').Value = '22.95'
$ws.Range('D32').Value = '7.23'
$ws.Range('E32').Value = '  +7.06%  '
$ws.Range('E34').Value = '  +3.07%  '
$ws.Range('E35').Value = '  +3.06%  '
$ws.Range('D36').Value = '161.66'
$ws.Range('E36').Value = '  -7.53%  '
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('D38').Value = '0.834'
$ws.Range('E38').Value = '  -2.49%  '
$ws.Range('D39').Value = '26.66'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').Value = '4.63'
$ws.Range('E40').Value = '  +7.93%  '
$ws.Range('E41').Value = '  +3.90%  '
$ws.Range('D42').Value = '2.60'
$ws.Range('E42').Value = '  +1.69%  '
$ws.Range('D43').Value = '41.37'
$ws.Range('E43').Value = '  +2.38%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '25.49'
$ws.Range('E44').Value = '  +2.98%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').Value = '347.20'
$ws.Range('E45').Value = '  +3.85%  '
$ws.Range('D46').Value = '0.0689'
$ws.Range('E46').Value = '  +2.40%  '
$ws.Range('D47').Value = '2.632.34'
$ws.Range('E47').Value = '  -4.71%  '
$ws.Range('E48').Value = '  +2.01%  '
$ws.Range('E49').Value = '  -0.64%  '
$ws.Range('E50').Value = '  +1.20%  '
$ws.Range('D51').Value = '6.21'
$ws.Range('E51').Value = '  +3.43%  '

# Restore default style on the cells we temporarily forced to text format,
# so no stray number-format styling is introduced versus the original file.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
